# Apply updates described by the diff to the "Inscricoes" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 17
$ws.Range("E17").Value = 89

# Row 18
$ws.Range("E18").Value = 89

# Row 29
$ws.Range("E29").Value = 15

# Row 36
$ws.Range("E36").Value = 74
$ws.Range("F36").Value = 28
$ws.Range("H36").Value = 28

# Row 43
$ws.Range("E43").Value = 19
$ws.Range("F43").Value = 13
$ws.Range("H43").Value = 13

# Row 44
$ws.Range("E44").Value = 22

# Row 47
$ws.Range("E47").Value = 48
$ws.Range("F47").Value = 28
$ws.Range("H47").Value = 28

# Row 62
$ws.Range("E62").Value = 33

# Row 63
$ws.Range("E63").Value = 21

# Row 69
$ws.Range("E69").Value = 15

# Row 76
$ws.Range("E76").Value = 36
$ws.Range("F76").Value = 12
$ws.Range("H76").Value = 12
